$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(633).Delete()
